# Add data for 2022-06-18 (sheet's "through" date moves from June 09 -> June 10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-06-10"

# 2. Update the header label for the most-recent-month column (column B).
$ws.Range("B1").Value = "June 2022 (through June 10)"

# 3. New carjacking counts added for 2022-06-18 across several neighborhoods/months.
$ws.Range("Z2").Value = 1    # Englewood, June 2018
$ws.Range("H5").Value = 1    # South Shore, June 2021
$ws.Range("AR7").Value = 1   # Grand Boulevard, June 2015
$ws.Range("H38").Value = 1   # Chatham, June 2021
$ws.Range("H68").Value = 2   # Kenwood, June 2021
$ws.Range("B86").Value = 1   # Rush & Division, June 2022 (through June 10)
$ws.Range("N96").Value = 1   # Wicker Park, June 2020

# 4. Existing counts incremented by the new data point.
$ws.Range("H14").Value = 3   # Austin, June 2021 (was 2)
$ws.Range("T94").Value = 2   # West Pullman, June 2019 (was 1)
